$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of product data (A162:B163)
$ws.Range("A162").Value = 12600935
$ws.Range("B162").Value = 105
$ws.Range("A163").Value = 12613955
$ws.Range("B163").Value = 105

# Match the formatting already used by column A on the data rows above
# (left/top aligned "General" number format) by copying the format from A161
$ws.Range("A161").Copy()
$ws.Range("A162:A163").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Scroll the view down near the newly added rows and select the last new cell,
# mirroring the author's navigation while adding the rows
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 139
$ws.Range("A163").Select()
